$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are stored as literal text (matches source formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.324.21"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.833.35"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "235.83"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").Value = "0.6044"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "0.06996"
$ws.Range("E8").Value = "  -4.68%  "
$ws.Range("D9").Value = "0.2779"
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("D10").Value = "23.64"
$ws.Range("E10").Value = "  -4.12%  "
$ws.Range("D11").Value = "0.07627"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").Value = "1.834.72"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D14").Value = "0.6330"
$ws.Range("E14").Value = "  -4.06%  "
$ws.Range("D15").Value = "0.000009904"
$ws.Range("E15").Value = "  -5.78%  "
$ws.Range("D16").Value = "78.07"
$ws.Range("E16").Value = "  -4.04%  "
$ws.Range("D17").Value = "29.002.23"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "5.622"
$ws.Range("E18").Value = "  -9.67%  "
$ws.Range("D19").Value = "218.34"
$ws.Range("E19").Value = "  -7.68%  "
$ws.Range("D21").Value = "11.62"
$ws.Range("E21").Value = "  -4.58%  "
$ws.Range("D22").Value = "6.934"
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("D23").Value = "1.006"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "156.42"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").Value = "8.012"
$ws.Range("E25").Value = "  -4.59%  "
$ws.Range("D26").Value = "0.1295"
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").Value = "16.58"
$ws.Range("E27").Value = "  -3.66%  "
$ws.Range("D28").Value = "0.06489"
$ws.Range("E28").Value = "  -5.18%  "
$ws.Range("D29").Value = "1.425"
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("D30").Value = "1.447"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").Value = "3.854"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("D32").Value = "3.808"
$ws.Range("E32").Value = "  -5.10%  "
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("E34").Value = "  -4.55%  "
$ws.Range("D35").Value = "0.6517"
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("D37").Value = "2.755"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "0.01757"
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("D39").Value = "6.591"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "1.147.97"
$ws.Range("E40").Value = "  -6.89%  "
$ws.Range("D41").Value = "0.8948"
$ws.Range("E41").Value = "  -5.11%  "
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "101.10"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "62.43"
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("D46").Value = "0.00000000113"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").Value = "1.626"
$ws.Range("E47").Value = "  -3.45%  "
$ws.Range("D48").Value = "8.532"
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("D49").Value = "0.4556"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "6.456"
$ws.Range("E50").Value = "  -6.09%  "
$ws.Range("D51").Value = "0.05503"
$ws.Range("E51").Value = "  -2.29%  "
